$d = $word.ActiveDocument

# Locate the paragraph that leaked the OpenAI API token ("API Token: sk-proj-...")
# and remove it completely (including its paragraph mark), per the commit
# "Removed all the secerets".
$paras = $d.Paragraphs
$count = $paras.Count
$targetIndex = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*sk-proj-*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ne $null) {
    $target = $paras.Item($targetIndex)
    $target.Range.Delete()

    # The secret paragraph merges away, leaving the following paragraph's
    # content; re-insert the leading space run that was originally the
    # start of the "API " run before the secret, so spacing/formatting of
    # the now-merged paragraph is preserved.
    $following = $paras.Item($targetIndex)
    $insertPoint = $d.Range($following.Range.Start, $following.Range.Start)
    $insertPoint.InsertBefore(" ")
}
